# Share Skill added, Edit and Delete in Manage Listing.
#
# Replaces the old SignUp / SignIn test-user data (Vidhya Venugan /
# vidhyav9@gmail.com / Ithika2015 / http://www.skillswap.pro/Home) with the
# new tester's data (Priyanka Meka / priyanka.mekha@gmail.com /
# Bollepalli88 / http://localhost:5000), and bumps the ShareSkill listing's
# start/end dates forward by three years.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# SignIn sheet - updated first so new shared strings are introduced in the
# same order the original authoring tool produced them in.
# ----------------------------------------------------------------------
$wsSignIn = $wb.Worksheets.Item("SignIn")

# A2 used to be a real web-address hyperlink; the new value is a plain
# (non-linked) localhost URL, so drop just that one hyperlink.
$hyperlinksToRemove = @()
foreach ($hl in $wsSignIn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hyperlinksToRemove += $hl
    }
}
foreach ($hl in $hyperlinksToRemove) {
    $hl.Delete()
}
$wsSignIn.Range("A2").Value = "http://localhost:5000"

# B2 keeps its mailto: hyperlink, just repointed at the new address.
foreach ($hl in $wsSignIn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.Address = "mailto:priyanka.mekha@gmail.com"
    }
}
$wsSignIn.Range("B2").Value = "priyanka.mekha@gmail.com"

# C2 keeps its existing mailto:Test@123 hyperlink/display untouched -
# only the visible password text changes.
$wsSignIn.Range("C2").Value = "Bollepalli88"

# ----------------------------------------------------------------------
# SignUp sheet
# ----------------------------------------------------------------------
$wsSignUp = $wb.Worksheets.Item("SignUp")

$wsSignUp.Range("A2").Value = "Priyanka"
$wsSignUp.Range("B2").Value = "Meka"

foreach ($hl in $wsSignUp.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$2') {
        $hl.Address = "mailto:priyanka.mekha@gmail.com"
    }
}
$wsSignUp.Range("C2").Value = "priyanka.mekha@gmail.com"

$wsSignUp.Range("D2").Value = "Bollepalli88"
$wsSignUp.Range("E2").Value = "Bollepalli88"

# ----------------------------------------------------------------------
# ShareSkill sheet - Startdate/Enddate move from 2019 to 2022 (same day of
# year). Go through a scratch cell + PasteSpecial(formats) round trip so
# the original quote-prefixed date style on H2 survives a plain .Value
# write (which would otherwise normalise it to the closest equivalent
# number-format style).
# ----------------------------------------------------------------------
$wsShareSkill = $wb.Worksheets.Item("ShareSkill")

$wsShareSkill.Range("H2").Copy() | Out-Null
$wsShareSkill.Range("Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsShareSkill.Range("H2").Value = 44663
$wsShareSkill.Range("Z1").Copy() | Out-Null
$wsShareSkill.Range("H2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsShareSkill.Range("Z1").Clear() | Out-Null

$wsShareSkill.Range("I2").Value = 44693

# ----------------------------------------------------------------------
# Selection / active-tab bookkeeping, matching the saved workbook state:
#   - ShareSkill: no longer the active tab, selection moves to C12
#   - SignUp: becomes the active tab, selection moves to C9
#   - SignIn: left completely alone (selection stays C2)
# Visit ShareSkill first, finish on SignUp so SignUp ends up active.
# ----------------------------------------------------------------------
[void]$wsShareSkill.Activate()
[void]$wsShareSkill.Range("C12").Select()

[void]$wsSignUp.Activate()
[void]$wsSignUp.Range("C9").Select()
